$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156, shifting existing rows 156-177 down to 157-178
$ws.Rows.Item(156).Insert()

$ws.Cells.Item(156, 1).Value = 11
$ws.Cells.Item(156, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(156, 3).Value = "Bíobío"
$ws.Cells.Item(156, 4).Value = 44918
$ws.Cells.Item(156, 5).Value = 8
$ws.Cells.Item(156, 6).Value = 100112032
$ws.Cells.Item(156, 7).Value = "Zapallo italiano"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 100
$ws.Cells.Item(156, 11).Value = 7000
$ws.Cells.Item(156, 12).Value = 7500
$ws.Cells.Item(156, 13).Value = 7250
$ws.Cells.Item(156, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(156, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(156, 16).Value = 145
$ws.Cells.Item(156, 17).Value = 50
$ws.Cells.Item(156, 18).Value = "Hortaliza"
